# hot fix - and some of the excel and scrolling done
#
# The "Tenant Passport ID Number" column (column C, with its header cell
# and all of its per-row formatting) is removed from the lease-agreement
# template sheet. Everything to the right of it (Start date, Finish, Rent
# amount, ... Deposit exchange rate) shifts one column to the left.
#
# Deleting the entire column (rather than just clearing its contents) is
# what also drops the now-orphaned "Tenant Passport ID Number" shared
# string, shrinks the used range from A1:L14 to A1:K14, and re-numbers the
# `spans`/`col` metadata - exactly matching a real "right-click column C ->
# Delete" in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("Tenant Passport ID Number") and shift D:L left to C:K.
$ws.Columns("C").Delete() | Out-Null

# After a column delete, Excel leaves the whole column that now occupies
# the deleted slot selected (top cell active) rather than the single cell
# that was selected before the edit.
$ws.Range("C1:C1048576").Select() | Out-Null
